# Generate Report for Archive
#
# The localization-status report is refreshed: the single row of sample
# data on each sheet moves from "Ready for handoff" to "In Translation"
# (Overview!E2:F2, zh-cn!C2, de-de!C2 all point at the same shared
# string). Once the text shrinks, the status columns that were sized to
# fit the old ("Ready for handoff") text are re-fit (narrowed) to the
# new, shorter text.

$wb = $excel.ActiveWorkbook

# Character-unit width that re-fits these "Status" columns to the new,
# shorter "In Translation" text (was sized for "Ready for handoff").
$newStatusColWidth = 12.5

# --- Overview sheet: Status columns for zh-cn (E) and de-de (F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F2").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
